# "Master Data Tables - Test Data / master-machine_spec" worksheet update.
# Replaces the sample USB-drive master-data rows with Dell Vostro laptop
# rows (English + Arabic translation), then tweaks the sheet's scroll/
# selection state and page setup to match the authored change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: English record (id 1001) ---
$ws.Range("B2").Value = "Vostro"
$ws.Range("C2").Value = "Dell"
$ws.Range("D2").Value = 3568
$ws.Range("E2").Value = "DKS"
$ws.Range("F2").Value = 1.454
$ws.Range("G2").Value = "To take enrollments"
$ws.Range("H2").Value = "eng"
$ws.Range("J2").Value = "superadmin"
$ws.Range("K2").Value = "now()"

# --- Row 3: Arabic record (id 1002) ---
$ws.Range("B3").Value = "ستر  "
$ws.Range("C3").Value = "دلّ  "
$ws.Range("D3").Value = 3568
$ws.Range("E3").Value = "DKS"
$ws.Range("F3").Value = 1.454
$ws.Range("G3").Value = "لأخذ التسجيلات"
$ws.Range("H3").Value = "ara"
$ws.Range("J3").Value = "superadmin"
$ws.Range("K3").Value = "now()"

# --- View: scroll so column C is leftmost, select everything below the
#     header/data rows (as the reviewer left it before saving) ---
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("A4:XFD1048576").Select()

# --- Page setup ---
$ws.PageSetup.PaperSize = 9      # xlPaperA4
$ws.PageSetup.Orientation = 1    # xlPortrait
